# Regional Availability Factor.xlsx - apply the "updated 4.0 files and mdl" edit.
#
# Semantic changes applied:
#   1. About!C1        : revision date bumped from 3/15/2024 to 3/28/2024 (serial 45366 -> 45379)
#   2. RAF-capacity!B24: capacity-credit RAF for "hydrogen combustion turbine" 0.3 -> 1
#   3. RAF-capacity!B25: capacity-credit RAF for "hydrogen combined cycle"   0.3 -> 1
#   4. The active/selected sheet moves from "RAF-generation" to "RAF-capacity",
#      matching activeTab going from 1 to 3 and tabSelected moving sheets.

$wb = $excel.ActiveWorkbook

# 1. Bump the "Source:" date on the About sheet.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# 2 & 3. Update the hydrogen plant capacity-credit RAF values.
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# 4. Make RAF-capacity the active sheet/tab, with B25 selected (last edited cell),
#    leaving RAF-generation no longer the selected tab.
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
